$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-31 Friday", 2) | Out-Null
$d.Content.Find.Execute("71×80=5680", $true, $false, $false, $false, $false, $true, 1, $false, "87×41=3567", 2) | Out-Null
$d.Content.Find.Execute("74×76=5624", $true, $false, $false, $false, $false, $true, 1, $false, "94×95=8930", 2) | Out-Null
$d.Content.Find.Execute("21×49=1029", $true, $false, $false, $false, $false, $true, 1, $false, "97×53=5141", 2) | Out-Null
$d.Content.Find.Execute("12×76=912", $true, $false, $false, $false, $false, $true, 1, $false, "40×16=640", 2) | Out-Null
$d.Content.Find.Execute("82×88=7216", $true, $false, $false, $false, $false, $true, 1, $false, "91×46=4186", 2) | Out-Null
$d.Content.Find.Execute("13×78=1014", $true, $false, $false, $false, $false, $true, 1, $false, "78×58=4524", 2) | Out-Null
$d.Content.Find.Execute("51×93=4743", $true, $false, $false, $false, $false, $true, 1, $false, "21×75=1575", 2) | Out-Null
$d.Content.Find.Execute("57×45=2565", $true, $false, $false, $false, $false, $true, 1, $false, "32×44=1408", 2) | Out-Null
$d.Content.Find.Execute("71×50=3550", $true, $false, $false, $false, $false, $true, 1, $false, "69×25=1725", 2) | Out-Null
$d.Content.Find.Execute("24×21=504", $true, $false, $false, $false, $false, $true, 1, $false, "54×76=4104", 2) | Out-Null
$d.Content.Find.Execute("85×85=7225", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=1344", 2) | Out-Null
$d.Content.Find.Execute("46×61=2806", $true, $false, $false, $false, $false, $true, 1, $false, "73×37=2701", 2) | Out-Null
$d.Content.Find.Execute("26×27=702", $true, $false, $false, $false, $false, $true, 1, $false, "33×99=3267", 2) | Out-Null
$d.Content.Find.Execute("51×74=3774", $true, $false, $false, $false, $false, $true, 1, $false, "42×78=3276", 2) | Out-Null
$d.Content.Find.Execute("19×74=1406", $true, $false, $false, $false, $false, $true, 1, $false, "76×71=5396", 2) | Out-Null
$d.Content.Find.Execute("81×70=5670", $true, $false, $false, $false, $false, $true, 1, $false, "25×88=2200", 2) | Out-Null
$d.Content.Find.Execute("58×78=4524", $true, $false, $false, $false, $false, $true, 1, $false, "96×33=3168", 2) | Out-Null
$d.Content.Find.Execute("52×45=2340", $true, $false, $false, $false, $false, $true, 1, $false, "73×24=1752", 2) | Out-Null
$d.Content.Find.Execute("13×81=1053", $true, $false, $false, $false, $false, $true, 1, $false, "52×73=3796", 2) | Out-Null
$d.Content.Find.Execute("57×44=2508", $true, $false, $false, $false, $false, $true, 1, $false, "85×65=5525", 2) | Out-Null
$d.Content.Find.Execute("22×33=726", $true, $false, $false, $false, $false, $true, 1, $false, "88×68=5984", 2) | Out-Null
$d.Content.Find.Execute("55×89=4895", $true, $false, $false, $false, $false, $true, 1, $false, "68×21=1428", 2) | Out-Null
$d.Content.Find.Execute("99×76=7524", $true, $false, $false, $false, $false, $true, 1, $false, "36×15=540", 2) | Out-Null
$d.Content.Find.Execute("68×64=4352", $true, $false, $false, $false, $false, $true, 1, $false, "41×97=3977", 2) | Out-Null
$d.Content.Find.Execute("65×73=4745", $true, $false, $false, $false, $false, $true, 1, $false, "34×31=1054", 2) | Out-Null
